$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.053.46'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.830.42'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.87'
$ws.Range("E5").Value = '  -0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6219'
$ws.Range("E6").Value = '  -6.60%  '

$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07536'
$ws.Range("E8").Value = '  +1.55%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.62'
$ws.Range("E9").Value = '  +6.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2918'
$ws.Range("E10").Value = '  -0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.73'
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07627'
$ws.Range("E12").Value = '  -1.74%  '

$ws.Range("D13").Value = '1.831.13'
$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6634'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.17'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009073'
$ws.Range("E17").Value = '  +8.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.995'
$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D19").Value = '29.045.45'

$ws.Range("D20").Value = '2.076.97'
$ws.Range("E20").Value = '  -0.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '224.57'
$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.32'
$ws.Range("E22").Value = '  -1.21%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.187'
$ws.Range("E24").Value = '  +0.49%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.38'
$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.408'
$ws.Range("E27").Value = '  -2.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1360'
$ws.Range("E28").Value = '  -3.23%  '

$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.498'
$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.209'
$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.030'
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.045'
$ws.Range("E33").Value = '  -1.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05215'
$ws.Range("E34").Value = '  -1.39%  '

$ws.Range("E35").Value = '  -1.86%  '

$ws.Range("E36").Value = '  +1.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7322'
$ws.Range("E37").Value = '  -1.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.645'
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").Value = '1.273.40'
$ws.Range("E39").Value = '  -2.04%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01784'
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.746'
$ws.Range("E41").Value = '  +0.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.320'
$ws.Range("E42").Value = '  +7.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8935'
$ws.Range("E43").Value = '  -4.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.96'
$ws.Range("E45").Value = '  -0.29%  '

$ws.Range("D46").Value = '1.975.50'
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5114'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("E49").Value = '  -0.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3960'
$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.666'
$ws.Range("E51").Value = '  -5.27%  '
